$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "42.183.17"
$ws.Range("E2").Value = "  -1.65%  "

$ws.Range("D3").Value = "2.244.89"
$ws.Range("E3").Value = "  -1.92%  "

$ws.Range("E4").Value = "  +0.04%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "247.18"
$ws.Range("E5").Value = "  -2.18%  "

$ws.Range("E6").Value = "  -2.96%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "74.18"
$ws.Range("E7").Value = "  -1.18%  "

$ws.Range("E8").Value = "  +0.08%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.618"
$ws.Range("E9").Value = "  -4.97%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "41.00"
$ws.Range("E10").Value = "  +4.77%  "

$ws.Range("E11").Value = "  -4.48%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "7.07"
$ws.Range("E12").Value = "  -6.16%  "

$ws.Range("E13").Value = "  -3.27%  "

$ws.Range("D14").Value = "2.579.64"
$ws.Range("E14").Value = "  -2.02%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "14.48"
$ws.Range("E15").Value = "  -4.41%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.851"
$ws.Range("E16").Value = "  -2.45%  "

$ws.Range("D17").Value = "2.237.13"
$ws.Range("E17").Value = "  -2.34%  "

$ws.Range("D18").Value = "42.021.56"
$ws.Range("E18").Value = "  -1.77%  "

$ws.Range("D19").Value = "0.0₃0976"
$ws.Range("E19").Value = "  -2.71%  "

$ws.Range("B20").Value = "Uniswap"
$ws.Range("C20").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.10"
$ws.Range("E20").Value = "  -2.11%  "

$ws.Range("B21").Value = "Litecoin"
$ws.Range("C21").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "71.72"
$ws.Range("E21").Value = "  -0.78%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "2.31"
$ws.Range("E22").Value = "  +7.04%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "231.10"
$ws.Range("E23").Value = "  -1.96%  "

$ws.Range("E24").Value = "  +0.07%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "11.07"
$ws.Range("E25").Value = "  -2.42%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "3.55"
$ws.Range("E26").Value = "  -8.38%  "

$ws.Range("B27").Value = "InternetComputer(DFINITY)"
$ws.Range("C27").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "7.57"
$ws.Range("E27").Value = "  +21.39%  "

$ws.Range("B28").Value = "PancakeSwap"
$ws.Range("C28").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.30"
$ws.Range("E28").Value = "  -4.42%  "

$ws.Range("E29").Value = "  -0.95%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "169.21"
$ws.Range("E30").Value = "  +1.12%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "20.65"
$ws.Range("E31").Value = "  -1.81%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.0825"
$ws.Range("E32").Value = "  -4.75%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.119"
$ws.Range("E33").Value = "  -6.69%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "29.96"
$ws.Range("E34").Value = "  -4.02%  "

$ws.Range("E35").Value = "  -2.94%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "4.53"
$ws.Range("E36").Value = "  -3.30%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "4.85"
$ws.Range("E37").Value = "  +1.03%  "

$ws.Range("E38").Value = "  -2.20%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "13.34"
$ws.Range("E39").Value = "  -2.15%  "

$ws.Range("E40").Value = "  -5.39%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "5.78"
$ws.Range("E41").Value = "  -2.97%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "108.75"
$ws.Range("E42").Value = "  +3.22%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.202"
$ws.Range("E43").Value = "  -3.56%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "60.59"
$ws.Range("E44").Value = "  -0.89%  "

$ws.Range("E45").Value = "  -5.20%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.1000"
$ws.Range("E46").Value = "  -1.40%  "

$ws.Range("E48").Value = "  -3.86%  "

$ws.Range("E49").Value = "  -1.76%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.25"
$ws.Range("E50").Value = "  -0.79%  "

$ws.Range("B51").Value = "HuobiToken"
$ws.Range("C51").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.70"
$ws.Range("E51").Value = "  -0.98%  "
